# Re-baseline refactor: append new milestone re-baseline rows below the
# existing "Project MM18 Notes" row (row 26), then move the active
# selection to J29 (mirrors the final cursor position recorded in the
# saved workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Re-baseline this quarter",
    "Re-baseline ALB/Programme milestones",
    "Re-baseline ALB/Programme cost",
    "Re-baseline ALB/Programme benefits",
    "Re-baseline IPDC milestones",
    "Re-baseline IPDC cost",
    "Re-baseline IPDC benefits",
    "Re-baseline HMT milestones",
    "Re-baseline HMT cost",
    "Re-baseline HMT benefits"
)

$startRow = 27
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $labels[$i]
    # Touch B:D so the row materializes with the same (default) style as
    # column A, matching the plain/unstyled look of the new rows.
    $ws.Cells.Item($r, 1).NumberFormat = "General"
    $ws.Cells.Item($r, 2).NumberFormat = "General"
    $ws.Cells.Item($r, 3).NumberFormat = "General"
    $ws.Cells.Item($r, 4).NumberFormat = "General"
}

$ws.Range("J29").Select()
